{"js": "// Fix typo: \"Is there association ...\" -> \"Is there an association ...\"\nconst body = context.document.body;\n\nconst q1 = body.search(\"Is there association between opioid prescribing rates and opioid use disorder?\", { matchCase: true });\nq1.load(\"text\");\nawait context.sync();\n\nif (q1.items.length > 0) {\n  q1.items[0].insertText(\n    \"Is there an association between opioid prescribing rates and opioid use disorder?\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Clean up the (unrelated-text-wise) run fragmentation in the \"Look for opioid use\n// disorder trends...\" bullet so it becomes a single run again - same text, just\n// re-written in place so the engine coalesces it into one run instead of four.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(\"Look for opioid use disorder trends among different populations\") === 0) {\n    para.getRange().insertText(para.text, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: fix typo \"Is there association ...\" -> \"Is there an association ...\"\n$find = $d.Content.Find\n$find.Execute(\n    \"Is there association between opioid prescribing rates and opioid use disorder?\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Is there an association between opioid prescribing rates and opioid use disorder?\",\n    2\n) | Out-Null\n\n# Change 2: the \"Look for opioid use disorder trends...\" bullet is split across several\n# runs with identical combined text; rewrite it in place so it collapses back to a\n# single run (no actual wording change).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"Look for opioid use disorder trends among different populations\")) {\n        $r = $d.Range($p.Range.Start, $p.Range.End - 1)\n        $newText = $r.Text\n        $r.Delete()\n        $r.InsertAfter($newText)\n    }\n}\n"}
